$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new rows (18 and 19) before the footer rows, pushing the
#    old rows 18+ (incl. the "NOMBRE DEL REPRESENTANTE LEGAL" block) down by 2.
$ws.Rows("18:19").Insert()

# 2) The old row 17 (bottom-bordered "last row" style) needs to move down to
#    row 19 stylistically, while row 17 itself becomes a "middle" row (same
#    look as row 16). Capture the old row17 formatting into row19 first...
$ws.Range("B17:J17").Copy($ws.Range("B19:J19"))

# ...then restyle row 17 to match row 16 (the "middle" row look).
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# 3) New row 18 also uses the "middle" row look.
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# 4) Fill in the period values for the four debt rows.
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2504"

# F/G already carry 56940/1423500 from the copy above; row 19 (the smaller,
# final period) has a different "Valor Mora" amount.
$ws.Range("F19").Value = 11388

# 5) Update the account summary total and period count.
$ws.Range("E11").Value = 182208
$ws.Range("F13").Value = 4
